# "added POM structuring to SKU"
# Rename the human-readable row labels in column A of Sheet1 to
# camelCase property names (Page-Object-Model style keys), leaving the
# data in column B untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value  = "itemHSNCode"
$ws.Range("A2").Value  = "skuName"
$ws.Range("A3").Value  = "localName"
$ws.Range("A4").Value  = "description"
$ws.Range("A5").Value  = "category"
$ws.Range("A6").Value  = "subCategory"
$ws.Range("A7").Value  = "gstRate"
$ws.Range("A8").Value  = "cessAmount"
$ws.Range("A9").Value  = "amountType"
$ws.Range("A10").Value = "brandName"
$ws.Range("A11").Value = "numberOfVariations"
$ws.Range("A12").Value = "variationName"
$ws.Range("A13").Value = "price"
$ws.Range("A14").Value = "valueText"
$ws.Range("A15").Value = "variationName"
$ws.Range("A16").Value = "price"
$ws.Range("A17").Value = "valueText"

# Move the active selection to A17, matching the saved view state.
$ws.Range("A17").Select()
